# Test Bounder Domain.xlsx - add Medicine + Caculator sheets, populate User sheet,
# and re-point the active sheet/selection the way the author's Excel session ended up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Populate the (currently empty) "User" sheet with the Add/Edit User spec.
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")

$wsUser.Columns.Item(1).ColumnWidth = 14.25
$wsUser.Columns.Item(2).ColumnWidth = 38.5

$wsUser.Range("A2:C2").Merge()
$wsUser.Range("A2").Value = "Add/Edit User"
$wsUser.Range("A2:C2").HorizontalAlignment = -4108
$wsUser.Range("A2:C2").VerticalAlignment = -4108

# Email
$wsUser.Range("A3").Value = "Email"
$wsUser.Range("B3").Value = "valid"
$wsUser.Range("C3").Value = "invalid"
$wsUser.Range("A4").Value = "format"
$wsUser.Range("B4").Value = "(characters)@(characters).(characters)"
$wsUser.Range("A5").Value = "length"
$wsUser.Range("B5").Value = "1-255"
$wsUser.Range("C5").Value = "0,256"
$loEmail = $wsUser.ListObjects.Add(1, $wsUser.Range("A3:C5"), $null, 1)
$loEmail.Name = "表2_1421"
$loEmail.TableStyle = "TableStyleLight9"

# Password
$wsUser.Range("A8").Value = "Password"
$wsUser.Range("B8").Value = "valid"
$wsUser.Range("C8").Value = "invalid"
$wsUser.Range("A9").Value = "character"
$wsUser.Range("B9").Value = "anything"
$wsUser.Range("A10").Value = "length"
$wsUser.Range("B10").Value = "1-255"
$wsUser.Range("C10").Value = "0,256"
$loPassword = $wsUser.ListObjects.Add(1, $wsUser.Range("A8:C10"), $null, 1)
$loPassword.Name = "表2_142122"
$loPassword.TableStyle = "TableStyleLight9"

# Password confirm
$wsUser.Range("A12").Value = "Password confirm"
$wsUser.Range("B12").Value = "valid"
$wsUser.Range("C12").Value = "invalid"
$wsUser.Range("A13").Value = "character"
$wsUser.Range("B13").Value = "anything"
$wsUser.Range("A14").Value = "length"
$wsUser.Range("B14").Value = "1-255"
$wsUser.Range("C14").Value = "0,256"
$loPasswordConfirm = $wsUser.ListObjects.Add(1, $wsUser.Range("A12:C14"), $null, 1)
$loPasswordConfirm.Name = "表2_14212223"
$loPasswordConfirm.TableStyle = "TableStyleLight9"

# FirstName
$wsUser.Range("A17").Value = "FirstName"
$wsUser.Range("B17").Value = "valid"
$wsUser.Range("C17").Value = "invalid"
$wsUser.Range("A18").Value = "character"
$wsUser.Range("B18").Value = "anything"
$wsUser.Range("A19").Value = "length"
$wsUser.Range("B19").Value = "1-255"
$wsUser.Range("C19").Value = "0,256"
$loFirstName = $wsUser.ListObjects.Add(1, $wsUser.Range("A17:C19"), $null, 1)
$loFirstName.Name = "表2_14212224"
$loFirstName.TableStyle = "TableStyleLight9"

# LastName
$wsUser.Range("A21").Value = "LastName"
$wsUser.Range("B21").Value = "valid"
$wsUser.Range("C21").Value = "invalid"
$wsUser.Range("A22").Value = "character"
$wsUser.Range("B22").Value = "anything"
$wsUser.Range("A23").Value = "length"
$wsUser.Range("B23").Value = "1-255"
$wsUser.Range("C23").Value = "0,256"
$loLastName = $wsUser.ListObjects.Add(1, $wsUser.Range("A21:C23"), $null, 1)
$loLastName.Name = "表2_1421222425"
$loLastName.TableStyle = "TableStyleLight9"

# Phone
$wsUser.Range("A26").Value = "Phone"
$wsUser.Range("B26").Value = "valid"
$wsUser.Range("C26").Value = "invalid"
$wsUser.Range("D26").Value = "列1"
$wsUser.Range("A27").Value = "character"
$wsUser.Range("B27").Value = "anything"
$wsUser.Range("D27").Value = "The reason that we didn't define it numeric type is there might be ""-"" in it."
$wsUser.Range("A28").Value = "length"
$wsUser.Range("B28").Value = "1-255"
$wsUser.Range("C28").Value = "0,256"
$loPhone = $wsUser.ListObjects.Add(1, $wsUser.Range("A26:D28"), $null, 1)
$loPhone.Name = "表2_142122242526"
$loPhone.TableStyle = "TableStyleLight9"

$wsUser.PageSetup.PaperSize = 9
$wsUser.PageSetup.Orientation = 1

$wsUser.Range("K27").Select()

# ---------------------------------------------------------------------------
# 2) Add "Medicine" sheet (Medicine Add/Update spec) after Visit.
# ---------------------------------------------------------------------------
$wsMedicine = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMedicine.Name = "Medicine"

$wsMedicine.Columns.Item(1).ColumnWidth = 47.75
$wsMedicine.Columns.Item(2).ColumnWidth = 32.25
$wsMedicine.Columns.Item(3).ColumnWidth = 48.25

$wsMedicine.Range("A2:C2").Merge()
$wsMedicine.Range("A2").Value = "Medicine Add/Update"
$wsMedicine.Range("A2:C2").HorizontalAlignment = -4108
$wsMedicine.Range("A2:C2").VerticalAlignment = -4108

# Medicine Name
$wsMedicine.Range("A3").Value = "Medicine Name"
$wsMedicine.Range("B3").Value = "valid"
$wsMedicine.Range("C3").Value = "invalid"
$wsMedicine.Range("A4").Value = "character"
$wsMedicine.Range("B4").Value = "a-z, A-Z, numbers, symbols"
$wsMedicine.Range("C4").Value = "others(blank)"
$wsMedicine.Range("A5").Value = "length"
$wsMedicine.Range("B5").Value = "1-255"
$wsMedicine.Range("C5").Value = "0,256"
$loMedName = $wsMedicine.ListObjects.Add(1, $wsMedicine.Range("A3:C5"), $null, 1)
$loMedName.Name = "表27"
$loMedName.TableStyle = "TableStyleLight9"

# min_dose
$wsMedicine.Range("A7").Value = "min_dose"
$wsMedicine.Range("B7").Value = "valid"
$wsMedicine.Range("C7").Value = "invalid"
$wsMedicine.Range("A8").Value = "character"
$wsMedicine.Range("B8").Value = "numbers(float or integers)"
$wsMedicine.Range("C8").Value = "others(letter or symbols)"
$wsMedicine.Range("A9").Value = "length"
$wsMedicine.Range("B9").Value = "(0,10^40) and (-10^40, 0)"
$wsMedicine.Range("C9").Value = "10^40+1,10^40"
$loMinDose = $wsMedicine.ListObjects.Add(1, $wsMedicine.Range("A7:C9"), $null, 1)
$loMinDose.Name = "表28"
$loMinDose.TableStyle = "TableStyleLight9"

# max_dose
$wsMedicine.Range("A11").Value = "max_dose"
$wsMedicine.Range("B11").Value = "valid"
$wsMedicine.Range("C11").Value = "invalid"
$wsMedicine.Range("A12").Value = "character"
$wsMedicine.Range("B12").Value = "numbers(float or integers)"
$wsMedicine.Range("C12").Value = "others(letter or symbols)"
$wsMedicine.Range("A13").Value = "length"
$wsMedicine.Range("B13").Value = "(0,10^40) and (-10^40, 0)"
$wsMedicine.Range("C13").Value = "10^40+1,10^40"
$loMaxDose = $wsMedicine.ListObjects.Add(1, $wsMedicine.Range("A11:C13"), $null, 1)
$loMaxDose.Name = "表29"
$loMaxDose.TableStyle = "TableStyleLight9"

# unit
$wsMedicine.Range("A15").Value = "unit"
$wsMedicine.Range("B15").Value = "valid"
$wsMedicine.Range("C15").Value = "invalid"
$wsMedicine.Range("A16").Value = "character"
$wsMedicine.Range("B16").Value = "number, letter,symbols"
$wsMedicine.Range("C16").Value = "others"
$wsMedicine.Range("A17").Value = "length"
$wsMedicine.Range("B17").Value = "0-8"
$wsMedicine.Range("C17").Value = "9 or more"
$loUnit = $wsMedicine.ListObjects.Add(1, $wsMedicine.Range("A15:C17"), $null, 1)
$loUnit.Name = "表30"
$loUnit.TableStyle = "TableStyleLight9"

# Hypo, weigtht, Renal_gu, Gi_sx,Chf, Cvd, Bone
$wsMedicine.Range("A19").Value = "Hypo, weigtht, Renal_gu, Gi_sx,Chf, Cvd, Bone"
$wsMedicine.Range("B19").Value = "valid"
$wsMedicine.Range("C19").Value = "invalid"
$wsMedicine.Range("A20").Value = "characters"
$wsMedicine.Range("B20").Value = "intergers 0-4"
$wsMedicine.Range("C20").Value = "others"
$wsMedicine.Range("A21").Value = "length"
$wsMedicine.Range("B21").Value = "0-4"
$wsMedicine.Range("C21").Value = "5,6 or others"
$loHypo = $wsMedicine.ListObjects.Add(1, $wsMedicine.Range("A19:C21"), $null, 1)
$loHypo.Name = "表31"
$loHypo.TableStyle = "TableStyleLight9"

$wsMedicine.Rows.Item(22).RowHeight = 21

$wsMedicine.Range("D4").Select()

# ---------------------------------------------------------------------------
# 3) Add empty "Caculator" sheet after Medicine; it ends up the active sheet.
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCalc.Name = "Caculator"

# ---------------------------------------------------------------------------
# 4) Re-point the Visit sheet's selection (it is no longer the active tab).
# ---------------------------------------------------------------------------
$wsVisit = $wb.Worksheets.Item("Visit")
$wsVisit.Range("A3:C5").Select()

# ---------------------------------------------------------------------------
# 5) Caculator ends up active, with F24 selected.
# ---------------------------------------------------------------------------
$wsCalc.Activate()
$wsCalc.Range("F24").Select()
